# Fruta / hortaliza, semanal
# Adds this week's new price observations for "Zapallo" (Camote / Paine)
# at the top of the "Vega Modelo de Temuco" data block, pushing the
# existing rows 941:964 down to 943:966.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 941 (shifts 941:964 -> 943:966
# and extends the used range to row 966).
$ws.Range("A941:A942").EntireRow.Insert()

# New row 941 - Camote
$ws.Range("A941").Value = 10
$ws.Range("B941").Value = "Vega Modelo de Temuco"
$ws.Range("C941").Value = "La Araucanía"
$ws.Range("D941").Value = 45239
$ws.Range("E941").Value = 9
$ws.Range("F941").Value = 100112045
$ws.Range("G941").Value = "Zapallo"
$ws.Range("H941").Value = "Camote"
$ws.Range("I941").Value = "1a (guarda)"
$ws.Range("J941").Value = 500
$ws.Range("K941").Value = 1600
$ws.Range("L941").Value = 1800
$ws.Range("M941").Value = 1640
$ws.Range("N941").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O941").Value = "Región del Maule"
$ws.Range("P941").Value = 1640
$ws.Range("Q941").Value = 1
$ws.Range("R941").Value = "Hortaliza"

# New row 942 - Paine
$ws.Range("A942").Value = 10
$ws.Range("B942").Value = "Vega Modelo de Temuco"
$ws.Range("C942").Value = "La Araucanía"
$ws.Range("D942").Value = 45239
$ws.Range("E942").Value = 9
$ws.Range("F942").Value = 100112045
$ws.Range("G942").Value = "Zapallo"
$ws.Range("H942").Value = "Paine"
$ws.Range("I942").Value = "1a (guarda)"
$ws.Range("J942").Value = 1000
$ws.Range("K942").Value = 1000
$ws.Range("L942").Value = 1000
$ws.Range("M942").Value = 1000
$ws.Range("N942").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O942").Value = "Región del Maule"
$ws.Range("P942").Value = 1000
$ws.Range("Q942").Value = 1
$ws.Range("R942").Value = "Hortaliza"
